# Update column F (dSF) values on Sheet1 to match repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    4  = -3
    5  = 2
    6  = 0
    7  = -2
    8  = 2
    10 = 9
    11 = -1
    12 = 6
    13 = 2
    14 = 1
    15 = -2
    16 = 2
    18 = -1
    19 = -7
    20 = 8
    21 = -2
    22 = -1
    24 = 8
    25 = -2
    27 = -3
    30 = 5
    31 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
